$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet 1: 展览 (Exhibition)
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# G2: number 88 -> text "不可售" (sold out / not for sale); F2 (30) unchanged
$ws1.Range("G2").Value = "不可售"

$ws1.Range("F3").Value = 646
$ws1.Range("F6").Value = 5646
$ws1.Range("F7").Value = 1586
$ws1.Range("F9").Value = 3207
$ws1.Range("F12").Value = 1338
$ws1.Range("F13").Value = 4474
$ws1.Range("F14").Value = 1070
$ws1.Range("F15").Value = 1693
$ws1.Range("F18").Value = 44
$ws1.Range("F20").Value = 166
$ws1.Range("F21").Value = 155
$ws1.Range("F22").Value = 1011
$ws1.Range("F24").Value = 81
$ws1.Range("F29").Value = 1110
$ws1.Range("F30").Value = 402
$ws1.Range("F32").Value = 195
$ws1.Range("F33").Value = 349
$ws1.Range("F34").Value = 258
$ws1.Range("F35").Value = 12
$ws1.Range("F36").Value = 1727
$ws1.Range("F37").Value = 2230
$ws1.Range("F38").Value = 1046
$ws1.Range("F42").Value = 351
$ws1.Range("F43").Value = 25
$ws1.Range("F45").Value = 24
$ws1.Range("F46").Value = 431
$ws1.Range("F47").Value = 380
$ws1.Range("F48").Value = 226

# ------------------------------------------------------------------
# Sheet 2: 演出 (Performance)
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("F22").Value = 10
$ws2.Range("F23").Value = 3
$ws2.Range("F24").Value = 2

# ------------------------------------------------------------------
# Sheet 3: 本地生活 (Local Life)
# ------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("F2").Value = 775

# ------------------------------------------------------------------
# Sheet 4: 全部类型 (All Types)
# ------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Range("F2").Value = 775

# Row 3 is fully replaced with a new event (old Aniidol listing -> DragonBand concert).
# B3 is a plain date-looking string; assigning it directly would get auto-coerced to an
# Excel date serial, so force-text it with a leading apostrophe and then reset the style
# back to Normal so no stray style index gets attached to the cell.
$ws4.Range("B3").Value = "'2024-04-30"
$ws4.Range("B3").Style = "Normal"
$ws4.Range("C3").Value = "杭州·DragonBand七龙珠限定Live纪念演出"
$ws4.Range("D3").Value = "中山南路77号尚城·利星1157 3F MAOLivehouse杭州"
$ws4.Range("E3").Value = "2024.04.30 20:30-04.30 22:30"
$ws4.Range("F3").Value = 4
$ws4.Range("G3").Value = 158
$ws4.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=84066"
$ws4.Range("I3").Value = "//i2.hdslb.com/bfs/openplatform/202403/j25rRJ8Q1711869134706.jpeg"

$ws4.Range("F5").Value = 5646
$ws4.Range("F6").Value = 1586
$ws4.Range("F9").Value = 3207
$ws4.Range("F10").Value = 1338
$ws4.Range("F11").Value = 4474
$ws4.Range("F12").Value = 1070
$ws4.Range("F16").Value = 44
$ws4.Range("F20").Value = 166
$ws4.Range("F21").Value = 155
$ws4.Range("F23").Value = 1011
$ws4.Range("F25").Value = 81
$ws4.Range("F30").Value = 1110
$ws4.Range("F31").Value = 402
$ws4.Range("F32").Value = 195
$ws4.Range("F33").Value = 258
$ws4.Range("F34").Value = 12
$ws4.Range("F35").Value = 1727
$ws4.Range("F37").Value = 1046
$ws4.Range("F42").Value = 351
$ws4.Range("F44").Value = 431
$ws4.Range("F45").Value = 380
$ws4.Range("F46").Value = 226
$ws4.Range("F47").Value = 2
